# Added user trial data for Subjects 1-6: raw measurements (rows 3-5 and
# row 13), the Average/Percent-Error formulas then recalculate automatically
# (clearing the #DIV/0! placeholders). Subject 3's B7 and Subject 5's E7
# use a flipped-sign variant of the percent-error formula, and Subject 4's
# D5 holds a malformed numeric entry ("0.0.0094") that Excel keeps as text.
# Finally, re-create each sheet's last-used selection/active-tab so the
# saved workbook lands on Subject 6 (the last sheet touched), matching the
# author's final editing position.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ----- Sheet 1 (Subject 1) -----
$ws = $sheets.Item(1)
$ws.Range("B3").Value = 0.1854
$ws.Range("C3").Value = 0.0413
$ws.Range("D3").Value = 0.0175
$ws.Range("E3").Value = 0.0094
$ws.Range("B4").Value = 0.1853
$ws.Range("C4").Value = 0.0392
$ws.Range("D4").Value = 0.0173
$ws.Range("E4").Value = 0.0094
$ws.Range("B5").Value = 0.1853
$ws.Range("C5").Value = 0.0391
$ws.Range("D5").Value = 0.0174
$ws.Range("E5").Value = 0.0092
$ws.Range("B13").Value = 0.212
$ws.Range("C13").Value = 0.046
$ws.Range("D13").Value = 0.0207
$ws.Range("E13").Value = 0.0094
$ws.Range("B6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C3:C5)"
$ws.Range("D6").Formula = "=AVERAGE(D3:D5)"
$ws.Range("E6").Formula = "=AVERAGE(E3:E5)"
$ws.Range("B7").Formula = "=(0.2-B6)*100/0.2"
$ws.Range("C7").Formula = "=(0.05-C6)*100/0.05"
$ws.Range("D7").Formula = "=(0.02-D6)*100/0.02"
$ws.Range("E7").Formula = "=(0.01-E6)*100/0.01"
[void]$ws.Range("E26").Select()

# ----- Sheet 2 (Subject 2) -----
$ws = $sheets.Item(2)
$ws.Range("B3").Value = 0.1822
$ws.Range("C3").Value = 0.0411
$ws.Range("D3").Value = 0.0183
$ws.Range("E3").Value = 0.0087
$ws.Range("B4").Value = 0.1786
$ws.Range("C4").Value = 0.0414
$ws.Range("D4").Value = 0.0187
$ws.Range("E4").Value = 0.0087
$ws.Range("B5").Value = 0.1782
$ws.Range("C5").Value = 0.0409
$ws.Range("D5").Value = 0.0184
$ws.Range("E5").Value = 0.0077
$ws.Range("B13").Value = 0.179
$ws.Range("C13").Value = 0.0455
$ws.Range("D13").Value = 0.0181
$ws.Range("E13").Value = 0.0111
$ws.Range("B6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C3:C5)"
$ws.Range("D6").Formula = "=AVERAGE(D3:D5)"
$ws.Range("E6").Formula = "=AVERAGE(E3:E5)"
$ws.Range("B7").Formula = "=(0.2-B6)*100/0.2"
$ws.Range("C7").Formula = "=(0.05-C6)*100/0.05"
$ws.Range("D7").Formula = "=(0.02-D6)*100/0.02"
$ws.Range("E7").Formula = "=(0.01-E6)*100/0.01"
[void]$ws.Range("G10").Select()

# ----- Sheet 3 (Subject 3) -----
$ws = $sheets.Item(3)
$ws.Range("B3").Value = 0.2761
$ws.Range("C3").Value = 0.0445
$ws.Range("D3").Value = 0.0061
$ws.Range("E3").Value = 0.0017
$ws.Range("B4").Value = 0.2791
$ws.Range("C4").Value = 0.0439
$ws.Range("D4").Value = 0.0065
$ws.Range("E4").Value = 0.0074
$ws.Range("B5").Value = 0.2743
$ws.Range("C5").Value = 0.0442
$ws.Range("D5").Value = 0.006
$ws.Range("E5").Value = 0.0078
$ws.Range("B13").Value = 0.1839
$ws.Range("C13").Value = 0.00478
$ws.Range("D13").Value = 0.0111
$ws.Range("E13").Value = 0.0109
$ws.Range("D14").Value = 0.0134
$ws.Range("D15").Value = 0.0311
$ws.Range("D16").Value = 0.0167
$ws.Range("D17").Value = 0.0187
$ws.Range("B6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C3:C5)"
$ws.Range("D6").Formula = "=AVERAGE(D3:D5)"
$ws.Range("E6").Formula = "=AVERAGE(E3:E5)"
$ws.Range("B7").Formula = "=(B6-0.2)*100/0.2"
$ws.Range("C7").Formula = "=(0.05-C6)*100/0.05"
$ws.Range("D7").Formula = "=(0.02-D6)*100/0.02"
$ws.Range("E7").Formula = "=(0.01-E6)*100/0.01"
[void]$ws.Range("F21").Select()

# ----- Sheet 4 (Subject 4) -----
$ws = $sheets.Item(4)
$ws.Range("B3").Value = 0.1847
$ws.Range("C3").Value = 0.0404
$ws.Range("D3").Value = 0.0083
$ws.Range("E3").Value = 0.0115
$ws.Range("B4").Value = 0.1843
$ws.Range("C4").Value = 0.0403
$ws.Range("D4").Value = 0.0163
$ws.Range("E4").Value = 0.0115
$ws.Range("B5").Value = 0.1845
$ws.Range("C5").Value = 0.0394
$ws.Range("D5").Value = "0.0.0094"
$ws.Range("E5").Value = 0.0117
$ws.Range("B13").Value = 0.1921
$ws.Range("C13").Value = 0.0488
$ws.Range("D13").Value = 0.0216
$ws.Range("E13").Value = 0.0094
$ws.Range("B6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C3:C5)"
$ws.Range("D6").Formula = "=AVERAGE(D3:D5)"
$ws.Range("E6").Formula = "=AVERAGE(E3:E5)"
$ws.Range("B7").Formula = "=(0.2-B6)*100/0.2"
$ws.Range("C7").Formula = "=(0.05-C6)*100/0.05"
$ws.Range("D7").Formula = "=(0.02-D6)*100/0.02"
$ws.Range("E7").Formula = "=(E6-0.01)*100/0.01"
[void]$ws.Range("E13").Select()

# ----- Sheet 5 (Subject 5) -----
$ws = $sheets.Item(5)
$ws.Range("B3").Value = 0.1915
$ws.Range("C3").Value = 0.0483
$ws.Range("D3").Value = 0.0201
$ws.Range("E3").Value = 0.0038
$ws.Range("B4").Value = 0.1913
$ws.Range("C4").Value = 0.051
$ws.Range("D4").Value = 0.02
$ws.Range("E4").Value = 0.0021
$ws.Range("B5").Value = 0.1913
$ws.Range("C5").Value = 0.0516
$ws.Range("D5").Value = 0.0198
$ws.Range("E5").Value = 0.0016
$ws.Range("B13").Value = 0.1922
$ws.Range("C13").Value = 0.0439
$ws.Range("D13").Value = 0.0201
$ws.Range("E13").Value = 0.0055
$ws.Range("E14").Value = 0.0106
$ws.Range("B6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C3:C5)"
$ws.Range("D6").Formula = "=AVERAGE(D3:D5)"
$ws.Range("E6").Formula = "=AVERAGE(E3:E5)"
$ws.Range("B7").Formula = "=(0.2-B6)*100/0.2"
$ws.Range("C7").Formula = "=(0.05-C6)*100/0.05"
$ws.Range("D7").Formula = "=(0.02-D6)*100/0.02"
$ws.Range("E7").Formula = "=(0.01-E6)*100/0.01"
[void]$ws.Range("E16").Select()

# ----- Sheet 6 (Subject 6) -----
$ws = $sheets.Item(6)
$ws.Range("B3").Value = 0.192
$ws.Range("C3").Value = 0.047
$ws.Range("D3").Value = 0.0178
$ws.Range("E3").Value = 0.0083
$ws.Range("B4").Value = 0.193
$ws.Range("C4").Value = 0.0458
$ws.Range("D4").Value = 0.0171
$ws.Range("E4").Value = 0.0081
$ws.Range("B5").Value = 0.1926
$ws.Range("C5").Value = 0.0458
$ws.Range("D5").Value = 0.017
$ws.Range("E5").Value = 0.0081
$ws.Range("B13").Value = 0.193
$ws.Range("C13").Value = 0.04661
$ws.Range("D13").Value = 0.0211
$ws.Range("E13").Value = 0.0095
$ws.Range("B6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C3:C5)"
$ws.Range("D6").Formula = "=AVERAGE(D3:D5)"
$ws.Range("E6").Formula = "=AVERAGE(E3:E5)"
$ws.Range("B7").Formula = "=(0.2-B6)*100/0.2"
$ws.Range("C7").Formula = "=(0.05-C6)*100/0.05"
$ws.Range("D7").Formula = "=(0.02-D6)*100/0.02"
$ws.Range("E7").Formula = "=(0.01-E6)*100/0.01"
[void]$ws.Range("E13").Select()
